$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = -0.1792682926829268
$ws.Range("G3").Value = -0.1792682926829268
$ws.Range("H2").Value = -0.3378048780487805
$ws.Range("H3").Value = -0.3378048780487805
$ws.Range("I2").Value = -0.2953982054892766
$ws.Range("I3").Value = -0.2953982054892766
$ws.Range("J2").Value = -0.2953982054892766
$ws.Range("J3").Value = -0.2953982054892766
$ws.Range("K2").Value = -38
$ws.Range("K3").Value = -38
$ws.Range("L2").Value = -0.4634146341463415
$ws.Range("L3").Value = -0.4634146341463415
$ws.Range("M2").Value = 0.154
$ws.Range("M3").Value = 0.154
$ws.Range("N2").Value = 0.0001745635910224439
$ws.Range("N3").Value = 0.0001745635910224439
$ws.Range("O2").Value = -0.004052631578947368
$ws.Range("O3").Value = -0.004052631578947368
$ws.Range("S2").Value = 0.154
$ws.Range("S3").Value = 0.154
$ws.Range("U2").Value = 81.40000000000001
$ws.Range("U3").Value = 81.40000000000001
$ws.Range("V2").Value = 0.09226932668329177
$ws.Range("V3").Value = 0.09226932668329177
$ws.Range("W2").Value = -0.8675799086757991
$ws.Range("W3").Value = -0.8675799086757991
$ws.Range("X2").Value = 0.1064708366164812
$ws.Range("X3").Value = 0.1064708366164812
$ws.Range("Y2").Value = -0.9740507452922804
$ws.Range("Y3").Value = -0.9740507452922804
$ws.Range("Z2").Value = 1.593369584966397
$ws.Range("Z3").Value = 1.593369584966397
$ws.Range("AA2").Value = -0.4706785160802672
$ws.Range("AA3").Value = -0.4706785160802672
$ws.Range("AB2").Value = 0.1034320223648314
$ws.Range("AB3").Value = 0.1034320223648314
$ws.Range("AC2").Value = -0.5741105384450986
$ws.Range("AC3").Value = -0.5741105384450986
$ws.Range("AD2").Value = 49.6
$ws.Range("AD3").Value = 49.6
$ws.Range("AE2").Value = 2.628264250603424
$ws.Range("AE3").Value = 2.628264250603424
$ws.Range("AF2").Value = 52.22826425060342
$ws.Range("AF3").Value = 52.22826425060342
$ws.Range("AG2").Value = -29.17173574939658
$ws.Range("AG3").Value = -29.17173574939658
$ws.Range("AH2").Value = 0.05589328389214517
$ws.Range("AH3").Value = 0.05589328389214517
$ws.Range("AI2").Value = 0.394985629937777
$ws.Range("AI3").Value = 0.394985629937777
$ws.Range("AJ2").Value = -0.03419785366083308
$ws.Range("AJ3").Value = -0.03419785366083308
$ws.Range("AK2").Value = -0.5739274433132008
$ws.Range("AK3").Value = -0.5739274433132008
$ws.Range("AL2").Value = 9.289999999999999
$ws.Range("AL3").Value = 9.289999999999999
$ws.Range("AM2").Value = 9.286999999999999
$ws.Range("AM3").Value = 9.286999999999999
$ws.Range("AN2").Value = -2.468030054236951
$ws.Range("AN3").Value = -2.468030054236951
$ws.Range("AO2").Value = -2.615715823466093
$ws.Range("AO3").Value = -2.615715823466093
$ws.Range("AP2").Value = 1.451546785559864
$ws.Range("AP3").Value = 1.451546785559864
$ws.Range("AQ2").Value = -2.616560783891462
$ws.Range("AQ3").Value = -2.616560783891462
